$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.587.92"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.882.73"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4737"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2888"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "99.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07827"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7570"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "1.881.58"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.233"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "30.570.75"
$ws.Range("E17").Value = "  -0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007518"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "2.128.04"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.342"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.416"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.157"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.905"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.325"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.503"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.240"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.179"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04842"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6982"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.770"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01903"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.869"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.297"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.967"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4234"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8401"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.989"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05787"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3948"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "
